$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.937.62"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.754.36"
$ws.Range("E3").Value = "  +0.87%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'579.53"
$ws.Range("E5").Value = "  -1.93%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'160.14"
$ws.Range("E6").Value = "  +4.75%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.39%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.40%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'5.86"
$ws.Range("E10").Value = "  -12.54%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +0.73%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.10%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.240.68"
$ws.Range("E13").Value = "  +0.69%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'27.07"
$ws.Range("E14").Value = "  +1.96%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "63.881.25"
$ws.Range("E15").Value = "  +0.01%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "'0.0000156"
$ws.Range("E16").Value = "  +2.02%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.757.58"
$ws.Range("E17").Value = "  +0.05%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +2.31%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.55%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'362.93"
$ws.Range("E20").Value = "  -0.60%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'6.90"
$ws.Range("E21").Value = "  -1.53%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.571"
$ws.Range("E22").Value = "  +6.15%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.46%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'66.42"
$ws.Range("E24").Value = "  +0.37%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  +2.98%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("E26").Value = "  +0.19%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.31%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0944"
$ws.Range("E28").Value = "  +4.26%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -1.40%  "

# Row 30 - Aptos
$ws.Range("E30").Value = "  +0.44%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  +5.15%  "

# Row 32 - Monero
$ws.Range("D32").Value = "'169.12"
$ws.Range("E32").Value = "  -2.18%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  +0.10%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "'20.61"
$ws.Range("E34").Value = "  +0.08%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "'5.03"
$ws.Range("E35").Value = "  +4.18%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +3.08%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  +2.39%  "

# Row 38 - SuiNetwork
$ws.Range("D38").Value = "'1.02"
$ws.Range("E38").Value = "  +0.59%  "

# Row 39 - Filecoin
$ws.Range("E39").Value = "  +0.12%  "

# Row 40 and 41 swap: RenderToken <-> Bittensor (with updated prices)
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'336.88"
$ws.Range("E40").Value = "  -2.92%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'6.15"
$ws.Range("E41").Value = "  +9.52%  "

# Row 42 - OKB
$ws.Range("D42").Value = "'39.53"
$ws.Range("E42").Value = "  +1.40%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "'22.21"
$ws.Range("E43").Value = "  +0.75%  "

# Row 44 - Hedera
$ws.Range("D44").Value = "'0.0602"
$ws.Range("E44").Value = "  +1.51%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").Value = "'22.08"
$ws.Range("E45").Value = "  -0.30%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +0.39%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  -1.21%  "

# Row 48 - Aave
$ws.Range("D48").Value = "'136.78"
$ws.Range("E48").Value = "  -4.61%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  +0.77%  "

# Row 50 - FirstDigitalUSD
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  +0.00%  "

# Row 51 - WhiteBITCoin
$ws.Range("E51").Value = "  +0.67%  "
